# Stock-report update: adjust item quantities (col F) and recompute their
# stock values (col G = Qty * Rate), then roll the deltas up into each
# company's "Sub Total:" cell (col B) and the two grand-total rows (724/725).
# A handful of adjacent line-item rows that share identical item names also
# swap their Item Code / Rate2 / Qty / Value (cols B, E, F, G) between the
# two rows - this mirrors an upstream re-sort of near-duplicate rows and
# does not change any totals (each pair's column sum is preserved).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F54").Value2 = 219
$ws.Range("G54").Value2 = 12285.9

$ws.Range("F57").Value2 = 125
$ws.Range("G57").Value2 = 11692.5

$ws.Range("F59").Value2 = 25
$ws.Range("G59").Value2 = 410.25

$ws.Range("B72").Value2 = 175991.87

$ws.Range("B132").Value2 = 64196
$ws.Range("F132").Value2 = 1
$ws.Range("G132").Value2 = 32143.58

$ws.Range("B133").Value2 = 65258
$ws.Range("F133").Value2 = 0
$ws.Range("G133").Value2 = 0

$ws.Range("B167").Value2 = 57756
$ws.Range("E167").Value2 = 79.37
$ws.Range("F167").Value2 = -100
$ws.Range("G167").Value2 = -6644

$ws.Range("B168").Value2 = 64350
$ws.Range("E168").Value2 = 70.63
$ws.Range("F168").Value2 = 2
$ws.Range("G168").Value2 = 132.88

$ws.Range("F184").Value2 = 92
$ws.Range("G184").Value2 = 8901.92

$ws.Range("B199").Value2 = 55961.14

$ws.Range("F216").Value2 = 96
$ws.Range("G216").Value2 = 5223.36

$ws.Range("F219").Value2 = 192
$ws.Range("G219").Value2 = 24322.56

$ws.Range("F220").Value2 = 87
$ws.Range("G220").Value2 = 9175.02

$ws.Range("F222").Value2 = 65
$ws.Range("G222").Value2 = 4829.5

$ws.Range("B224").Value2 = 65560.17999999999

$ws.Range("F226").Value2 = 37
$ws.Range("G226").Value2 = 2337.66

$ws.Range("F228").Value2 = 333
$ws.Range("G228").Value2 = 6160.5

$ws.Range("F229").Value2 = 17
$ws.Range("G229").Value2 = 364.31

$ws.Range("B235").Value2 = 14142.69

$ws.Range("F293").Value2 = 37
$ws.Range("G293").Value2 = 2025.38

$ws.Range("F297").Value2 = 17
$ws.Range("G297").Value2 = 1456.56

$ws.Range("F299").Value2 = 1
$ws.Range("G299").Value2 = 87.7

$ws.Range("B301").Value2 = 96643.49000000001

$ws.Range("B303").Value2 = 61610
$ws.Range("E303").Value2 = 122.71
$ws.Range("F303").Value2 = -58
$ws.Range("G303").Value2 = -5957.18

$ws.Range("B304").Value2 = 63565
$ws.Range("E304").Value2 = 109.19
$ws.Range("F304").Value2 = 60
$ws.Range("G304").Value2 = 6162.6

$ws.Range("B312").Value2 = 63531
$ws.Range("E312").Value2 = 152.53
$ws.Range("F312").Value2 = 24
$ws.Range("G312").Value2 = 3443.52

$ws.Range("B313").Value2 = 57802
$ws.Range("E313").Value2 = 162.71
$ws.Range("F313").Value2 = -79
$ws.Range("G313").Value2 = -11334.92

$ws.Range("F357").Value2 = 196
$ws.Range("G357").Value2 = 28349.44

$ws.Range("B362").Value2 = 72366.46000000001

$ws.Range("F372").Value2 = 43
$ws.Range("G372").Value2 = 2379.19

$ws.Range("B378").Value2 = 49177.76

$ws.Range("F391").Value2 = 29
$ws.Range("G391").Value2 = 862.46

$ws.Range("F393").Value2 = 361
$ws.Range("G393").Value2 = 34872.6

$ws.Range("B395").Value2 = 50932.92

$ws.Range("F402").Value2 = 114
$ws.Range("G402").Value2 = 2904.72

$ws.Range("F403").Value2 = 57
$ws.Range("G403").Value2 = 2061.69

$ws.Range("F404").Value2 = 82
$ws.Range("G404").Value2 = 4602.66

$ws.Range("F408").Value2 = 17
$ws.Range("G408").Value2 = 583.27

$ws.Range("F416").Value2 = 65
$ws.Range("G416").Value2 = 15670.2

$ws.Range("F419").Value2 = 67
$ws.Range("G419").Value2 = 3857.19

$ws.Range("F421").Value2 = 51
$ws.Range("G421").Value2 = 2789.7

$ws.Range("F422").Value2 = 48
$ws.Range("G422").Value2 = 1409.28

$ws.Range("B423").Value2 = 156419.04

$ws.Range("B485").Value2 = 64810
$ws.Range("E485").Value2 = 291.22
$ws.Range("F485").Value2 = 0
$ws.Range("G485").Value2 = 0

$ws.Range("B486").Value2 = 53319
$ws.Range("E486").Value2 = 310.64
$ws.Range("F486").Value2 = -6
$ws.Range("G486").Value2 = -1643.52

$ws.Range("B502").Value2 = 60025
$ws.Range("E502").Value2 = 37.22
$ws.Range("F502").Value2 = -98
$ws.Range("G502").Value2 = -3217.34

$ws.Range("B503").Value2 = 64833
$ws.Range("E503").Value2 = 34.9
$ws.Range("F503").Value2 = 88
$ws.Range("G503").Value2 = 2889.04

$ws.Range("B512").Value2 = 64830
$ws.Range("E512").Value2 = 34.9
$ws.Range("F512").Value2 = 83
$ws.Range("G512").Value2 = 2724.89

$ws.Range("B513").Value2 = 60022
$ws.Range("E513").Value2 = 37.22
$ws.Range("F513").Value2 = -113
$ws.Range("G513").Value2 = -3709.79

$ws.Range("F527").Value2 = 64
$ws.Range("G527").Value2 = 1753.6

$ws.Range("B531").Value2 = 108496.52

$ws.Range("F550").Value2 = 27
$ws.Range("G550").Value2 = 1671.3

$ws.Range("B562").Value2 = 37341.93

$ws.Range("F565").Value2 = 13
$ws.Range("G565").Value2 = 1351.22

$ws.Range("B567").Value2 = 19092.17

$ws.Range("F571").Value2 = 9
$ws.Range("G571").Value2 = 2528.55

$ws.Range("B579").Value2 = 13545.68

$ws.Range("F618").Value2 = 220
$ws.Range("G618").Value2 = 33090.2

$ws.Range("F620").Value2 = 73
$ws.Range("G620").Value2 = 10590.84

$ws.Range("F621").Value2 = 93
$ws.Range("G621").Value2 = 14384.31

$ws.Range("B634").Value2 = 193397.85

$ws.Range("F680").Value2 = 418
$ws.Range("G680").Value2 = 68179.98

$ws.Range("B686").Value2 = 69192.53

$ws.Range("F704").Value2 = 15
$ws.Range("G704").Value2 = 2542.95

$ws.Range("F711").Value2 = 1
$ws.Range("G711").Value2 = 679.58

$ws.Range("B719").Value2 = 57043.49

$ws.Range("B724").Value2 = 2319423.03

$ws.Range("B725").Value2 = 2319423.03
